$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style of existing header cell (G1) to new header cell (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
